# Included FLD_DocumentRegistry_New_Transmittals test details
# (Fulcrum_FluidTX_Trunk/src/com/proj/config/Config.xlsx)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # DataFetchFlag
$ws2 = $wb.Worksheets.Item(2)   # DataFetchXL

# --- Sheet "DataFetchFlag": new row 11 ---
$ws1.Range("A11").Value = 'FLD_DocumentRegistry_New_Transmittals'
$ws1.Range("B11").Value = 'XL'

# Extend the data-validation list from B2:B10 to B2:B11
$ws1.Range("B2:B10").Validation.Delete()
$ws1.Range("B2:B11").Validation.Add(3, 1, 1, """XL,DB""")

# --- Sheet "DataFetchXL": new row 11 ---
$ws2.Range("A11").Value = 'FLD_DocumentRegistry_New_Transmittals'
$ws2.Range("B11").Value = '\\src\\com\\proj\\suiteDOCS\\testdata\\DocumentRegistryTestData-Newtransmittal.xlsx'
$ws2.Range("C11").Value = 'Transmittals_New'

[void]$ws2.Hyperlinks.Add($ws2.Range("B11"), "file:///\\src\com\proj\suiteDOCS\testdata\DocumentRegistryTestData-Newtransmittal.xlsx")
$ws2.Range("B11").Style = "Hyperlink"

# --- Selections (match last-edited cell per sheet) ---
[void]$ws2.Range("C11").Select()
[void]$ws1.Range("A14").Select()
